$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings remain stored as text (matches source data format)
$textCells = @("D5", "D8", "D10", "D13", "D15", "D17", "D18", "D20", "D21", "D22", "D23", "D24", "D27", "D30", "D31", "D33", "D37", "D40", "D43", "D44", "D45", "D46", "D49")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply updated values
$ws.Range("D2").Value = '34.588.60'
$ws.Range("E2").Value = '  +1.30%  '
$ws.Range("D3").Value = '1.798.62'
$ws.Range("E3").Value = '  +1.07%  '
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("D5").Value = '227.29'
$ws.Range("E5").Value = '  +0.52%  '
$ws.Range("E6").Value = '  +1.93%  '
$ws.Range("E7").Value = '  -0.12%  '
$ws.Range("D8").Value = '32.83'
$ws.Range("E8").Value = '  +3.53%  '
$ws.Range("E9").Value = '  +1.38%  '
$ws.Range("D10").Value = '0.0695'
$ws.Range("E10").Value = '  +0.47%  '
$ws.Range("E11").Value = '  +0.46%  '
$ws.Range("D12").Value = '2.060.65'
$ws.Range("E12").Value = '  +1.25%  '
$ws.Range("D13").Value = '11.13'
$ws.Range("E13").Value = '  +1.83%  '
$ws.Range("D14").Value = '1.786.91'
$ws.Range("E14").Value = '  -0.07%  '
$ws.Range("D15").Value = '0.639'
$ws.Range("E15").Value = '  +2.85%  '
$ws.Range("D16").Value = '34.596.57'
$ws.Range("E16").Value = '  +1.45%  '
$ws.Range("D17").Value = '4.32'
$ws.Range("E17").Value = '  +3.11%  '
$ws.Range("D18").Value = '68.89'
$ws.Range("E18").Value = '  +1.46%  '
$ws.Range("E19").Value = '  +0.61%  '
$ws.Range("D20").Value = '246.84'
$ws.Range("E20").Value = '  +0.57%  '
$ws.Range("D21").Value = '11.38'
$ws.Range("E21").Value = '  +3.59%  '
$ws.Range("D22").Value = '1.00'
$ws.Range("E22").Value = '  -0.15%  '
$ws.Range("D23").Value = '4.17'
$ws.Range("E23").Value = '  +1.91%  '
$ws.Range("D24").Value = '173.41'
$ws.Range("E24").Value = '  +6.72%  '
$ws.Range("E25").Value = '  +1.75%  '
$ws.Range("E26").Value = '  +1.67%  '
$ws.Range("D27").Value = '16.65'
$ws.Range("E27").Value = '  +2.22%  '
$ws.Range("E28").Value = '  +1.86%  '
$ws.Range("E29").Value = '  -0.12%  '
$ws.Range("D30").Value = '4.04'
$ws.Range("E30").Value = '  +8.60%  '
$ws.Range("D31").Value = '0.0526'
$ws.Range("E31").Value = '  +0.96%  '
$ws.Range("E32").Value = '  +0.55%  '
$ws.Range("D33").Value = '3.80'
$ws.Range("E33").Value = '  +1.37%  '
$ws.Range("E34").Value = '  +2.01%  '
$ws.Range("D35").Value = '1.430.69'
$ws.Range("E35").Value = '  -0.65%  '
$ws.Range("E36").Value = '  +6.80%  '
$ws.Range("D37").Value = '0.678'
$ws.Range("E37").Value = '  +2.56%  '
$ws.Range("E38").Value = '  +2.59%  '
$ws.Range("E39").Value = '  +0.43%  '
$ws.Range("D40").Value = '84.77'
$ws.Range("E40").Value = '  +5.65%  '
$ws.Range("E41").Value = '  +2.97%  '
$ws.Range("E42").Value = '  +1.53%  '
$ws.Range("D43").Value = '2.75'
$ws.Range("E43").Value = '  +3.28%  '
$ws.Range("D44").Value = '13.87'
$ws.Range("E44").Value = '  +3.17%  '
$ws.Range("D45").Value = '0.0525'
$ws.Range("E45").Value = '  +2.75%  '
$ws.Range("D46").Value = '6.10'
$ws.Range("E46").Value = '  +0.49%  '
$ws.Range("E47").Value = '  +0.66%  '
$ws.Range("D48").Value = '1.960.92'
$ws.Range("E48").Value = '  +1.25%  '
$ws.Range("D49").Value = '105.16'
$ws.Range("E49").Value = '  +0.97%  '
$ws.Range("E50").Value = '  -0.16%  '
$ws.Range("E51").Value = '  -5.18%  '
